$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[53.950525608403936, 73.08834038903932]"
$ws.Range("U2").Value = "[43.72772818883496, 56.518398959786154]"
$ws.Range("M3").Value = "[53.55876814085056, 73.27247175096153]"
$ws.Range("N3").Value = [double]"2.220446049250313e-16"
$ws.Range("O3").Value = [double]"2.220446049250313e-16"
$ws.Range("U3").Value = "[41.683090009276334, 54.51522887129203]"
$ws.Range("M4").Value = "[52.88993708917917, 73.97381740775116]"
$ws.Range("N4").Value = [double]"8.881784197001252e-16"
$ws.Range("O4").Value = [double]"8.881784197001252e-16"
$ws.Range("U4").Value = "[43.78650224801575, 56.70053648043104]"
$ws.Range("M5").Value = "[52.54102138785514, 74.27099895352734]"
$ws.Range("N5").Value = [double]"2.664535259100376e-15"
$ws.Range("O5").Value = [double]"2.664535259100376e-15"
$ws.Range("U5").Value = "[43.85550158494461, 56.86358404626947]"
$ws.Range("M6").Value = "[52.345351224601465, 74.34733903327171]"
$ws.Range("N6").Value = [double]"4.218847493575595e-15"
$ws.Range("O6").Value = [double]"4.218847493575595e-15"
$ws.Range("U6").Value = "[41.808458799724235, 54.85321326161416]"
$ws.Range("M7").Value = "[51.255538556532755, 75.61326789393692]"
$ws.Range("N7").Value = [double]"1.13464793116691e-13"
$ws.Range("O7").Value = [double]"1.13464793116691e-13"
$ws.Range("U7").Value = "[43.9325828297452, 57.00595039132513]"
$ws.Range("M8").Value = "[51.0458732266313, 76.47210161757953]"
$ws.Range("N8").Value = [double]"3.790301406070284e-13"
$ws.Range("O8").Value = [double]"3.790301406070284e-13"
$ws.Range("U8").Value = "[43.191996676020565, 56.17168365105837]"
$ws.Range("M9").Value = "[50.99584731146365, 76.52212753274718]"
$ws.Range("N9").Value = [double]"4.289901767151605e-13"
$ws.Range("O9").Value = [double]"4.289901767151605e-13"
$ws.Range("U9").Value = "[43.19553027097845, 56.16815005610049]"
$ws.Range("M10").Value = "[51.88486223598158, 75.7254202271834]"
$ws.Range("N10").Value = [double]"4.685141163918161e-14"
$ws.Range("O10").Value = [double]"4.685141163918161e-14"
$ws.Range("U10").Value = "[43.230800245503815, 56.22171388397314]"
$ws.Range("M11").Value = "[51.348899223935916, 76.26138323922908]"
$ws.Range("N11").Value = [double]"1.938449400995523e-13"
$ws.Range("O11").Value = [double]"1.938449400995523e-13"
$ws.Range("U11").Value = "[43.235908081483416, 56.21660604799354]"
$ws.Range("M12").Value = "[51.162284349266486, 76.4479981138985]"
$ws.Range("N12").Value = [double]"3.110844914999689e-13"
$ws.Range("O12").Value = [double]"3.110844914999689e-13"
$ws.Range("Q12").Value = "[0.7736053981812692, 1.1509738850989635]"
$ws.Range("R12").Value = [double]"2.229327833447314e-13"
$ws.Range("S12").Value = [double]"2.229327833447314e-13"
$ws.Range("U12").Value = "[43.23804998994761, 56.21446413952934]"
$ws.Range("Y12").Value = 19.46474474474498
$ws.Range("Z12").Value = 20.89597597597623
